$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-12T03:01:46.787703+00:00"
$ws.Range("K3").Value = "2025-12-12T03:01:46.787736+00:00"
$ws.Range("K4").Value = "2025-12-12T03:01:46.787756+00:00"
$ws.Range("K5").Value = "2025-12-12T03:01:49.166917+00:00"
$ws.Range("K6").Value = "2025-12-12T03:01:49.166945+00:00"
$ws.Range("K7").Value = "2025-12-12T03:01:49.166962+00:00"
$ws.Range("K8").Value = "2025-12-12T03:01:51.998227+00:00"
$ws.Range("K9").Value = "2025-12-12T03:01:54.763394+00:00"
$ws.Range("K10").Value = "2025-12-12T03:01:57.589429+00:00"
$ws.Range("K11").Value = "2025-12-12T03:02:00.249649+00:00"
$ws.Range("K12").Value = "2025-12-12T03:02:05.310502+00:00"
$ws.Range("K13").Value = "2025-12-12T03:02:05.310537+00:00"
$ws.Range("K14").Value = "2025-12-12T03:02:08.255500+00:00"
$ws.Range("K15").Value = "2025-12-12T03:02:11.045341+00:00"
$ws.Range("K16").Value = "2025-12-12T03:02:13.820311+00:00"
$ws.Range("K17").Value = "2025-12-12T03:02:16.148308+00:00"
$ws.Range("K18").Value = "2025-12-12T03:02:16.148341+00:00"
$ws.Range("K19").Value = "2025-12-12T03:02:21.189985+00:00"
$ws.Range("K20").Value = "2025-12-12T03:02:21.190016+00:00"
$ws.Range("K21").Value = "2025-12-12T03:02:21.190034+00:00"
$ws.Range("K22").Value = "2025-12-12T03:02:23.927787+00:00"
$ws.Range("K23").Value = "2025-12-12T03:02:23.927815+00:00"
$ws.Range("K24").Value = "2025-12-12T03:02:23.927834+00:00"
$ws.Range("K25").Value = "2025-12-12T03:02:23.927851+00:00"
$ws.Range("K26").Value = "2025-12-12T03:02:23.927868+00:00"
$ws.Range("K27").Value = "2025-12-12T03:02:26.631590+00:00"
$ws.Range("K28").Value = "2025-12-12T03:02:26.631623+00:00"
$ws.Range("K29").Value = "2025-12-12T03:02:26.631643+00:00"
$ws.Range("K30").Value = "2025-12-12T03:02:29.622752+00:00"
$ws.Range("K31").Value = "2025-12-12T03:02:29.622781+00:00"
$ws.Range("K32").Value = "2025-12-12T03:02:36.037309+00:00"
$ws.Range("K33").Value = "2025-12-12T03:02:38.271901+00:00"
$ws.Range("K34").Value = "2025-12-12T03:02:38.271929+00:00"
$ws.Range("K35").Value = "2025-12-12T03:02:38.271946+00:00"
$ws.Range("K36").Value = "2025-12-12T03:02:40.997977+00:00"
$ws.Range("K37").Value = "2025-12-12T03:02:40.998007+00:00"
$ws.Range("K38").Value = "2025-12-12T03:02:43.965826+00:00"
$ws.Range("K39").Value = "2025-12-12T03:02:43.965870+00:00"
$ws.Range("K40").Value = "2025-12-12T03:02:46.763880+00:00"
$ws.Range("K41").Value = "2025-12-12T03:02:46.763911+00:00"
